# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# Replaces the 15 worker/period rows (rows 16-30) of the "Estado de Cuenta"
# table on Hoja1 with the updated roster + refreshed Valor Mora / Salario
# Basico figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @('CC', '45489600',   'SANDRA DE JESUS VARGAS BARRIOS',   '1805', 15624, 0),
    @('CC', '45489600',   'SANDRA DE JESUS VARGAS BARRIOS',   '1804', 15624, 0),
    @('CC', '32938139',   'CATHERINE MOSCOTE GELIS',          '1708', 7377,  737717),
    @('CC', '45557121',   'INGRID MARIA POLO ALIAN',          '1805', 15624, 781242),
    @('CC', '45557121',   'INGRID MARIA POLO ALIAN',          '1804', 15624, 781242),
    @('CC', '73204142',   'JEFFERSON VIAÑA REYES',            '1805', 15624, 781242),
    @('CC', '73204142',   'JEFFERSON VIAÑA REYES',            '1804', 15624, 781242),
    @('CC', '22798615',   'EDITH CORPAS DE AVILA',            '1708', 7377,  737717),
    @('CC', '1047476121', 'EDGAR DAVID GRÜNEWALD FLOREZ',     '1805', 15885, 794250),
    @('CC', '1047476121', 'EDGAR DAVID GRÜNEWALD FLOREZ',     '1804', 15885, 794250),
    @('CC', '45553335',   'LORENA PATRICIA RIPOLL FIGUEROA',  '1805', 15624, 781242),
    @('CC', '45553335',   'LORENA PATRICIA RIPOLL FIGUEROA',  '1804', 15624, 781242),
    @('CC', '1047464421', 'LAURA EUGENIA THERAN VASQUEZ',     '1805', 15885, 794250),
    @('CC', '1047464421', 'LAURA EUGENIA THERAN VASQUEZ',     '1804', 15885, 794250),
    @('CC', '73182307',   'YIMMY PATERNINA ARAUJO',           '1805', 15624, 781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $data[$i][0]
    $ws.Cells.Item($r, 3).Value = $data[$i][1]
    $ws.Cells.Item($r, 4).Value = $data[$i][2]
    $ws.Cells.Item($r, 5).Value = $data[$i][3]
    $ws.Cells.Item($r, 6).Value = $data[$i][4]
    $ws.Cells.Item($r, 7).Value = $data[$i][5]
}
